# Apply updated cryptocurrency price/volume data to Sheet1 (rows 2-51).
# Source: commit "Updated cryptos list on Sun Jul 30 09:50:59 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A style reference (no explicit number format / xf index 0) used to restore
# a cells style after temporarily forcing text format so that purely-numeric-
# looking strings (e.g. "1.177") are not auto-converted into numbers by Excel.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "29.309.46"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.875.18"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7128"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.38"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08027"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  +3.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3149"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08215"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("D12").Value = "1.877.36"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.97"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  +4.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.246"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7114"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.422"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +6.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008510"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +4.35%  "
$ws.Range("D18").Value = "29.313.29"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.68"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "2.134.51"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.24"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.763"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.037"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.40"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.502"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.403"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.300"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05368"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.180"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  -8.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.936"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7660"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  +3.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.178"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.691"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01875"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").Value = "1.260.16"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.435"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9125"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +3.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "112.40"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  +3.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "73.85"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  +2.10%  "
$ws.Range("E45").Value = "  +9.73%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "2.032.20"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5223"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.796"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.478"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4347"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  +1.06%  "
